$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Loan_Id" column right after "Sr_No" (old column B, Membership_Id, shifts to C)
$ws.Columns("B").Insert() | Out-Null
$ws.Range("B1").Value = "Loan_Id"

# Insert a new "Late_Fees" column right after "Membership_Id" (old numeric columns shift right by one more)
$ws.Columns("D").Insert() | Out-Null
$ws.Range("D1").Value = "Late_Fees"

# The old "Late_Fees(if_applicable)" header column has now shifted all the way to AI;
# remove it since the data it held is replaced by the new Late_Fees column above.
$ws.Columns("AI").Clear() | Out-Null

# Restore the selection to match the last saved view state.
$ws.Range("AJ2").Select() | Out-Null
